$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("status_label") that holds a string version of the
# existing "statut" emoji column A. This shifts the old B..I columns
# (NCTId, eudraCT, completion_year, clinical_trial_title, acronym,
# results_1y, results_3y, results) one place to the right, to C..J.
$ws.Columns("B").Insert()

# Header
$ws.Range("B1").Value = "status_label"

# Map each row's status emoji (column A) to its string label in new column B
$ws.Range("B2").Value = "orange"
$ws.Range("B3").Value = "rouge"
$ws.Range("B4").Value = "rouge"
$ws.Range("B5").Value = "rouge"
$ws.Range("B6").Value = "orange"
$ws.Range("B7").Value = "rouge"
